$wb = $excel.ActiveWorkbook

# Sheet 1: summ13356575 -> summ50993138
$ws = $wb.Worksheets.Item(1)
$ws.Name = 'summ50993138'
$arr1 = New-Object 'object[,]' 19,3
$arr1[0,0] = 'Intercept'
$arr1[0,1] = 42917.38032555665
$arr1[0,2] = 0.001100266509406315
$arr1[1,0] = 'Education[T.Primary/None]'
$arr1[1,1] = -1036.343153769936
$arr1[1,2] = 0.5863744054298416
$arr1[2,0] = 'Education[T.Secondary]'
$arr1[2,1] = -2254.570789624536
$arr1[2,2] = 0.1173677733020186
$arr1[3,0] = 'Education[T.University]'
$arr1[3,1] = -230.5386554266832
$arr1[3,2] = 0.6805518403813388
$arr1[4,0] = 'Season[T.Spring]'
$arr1[4,1] = -66.06597488323854
$arr1[4,2] = 0.9162514962012712
$arr1[5,0] = 'Season[T.Summer]'
$arr1[5,1] = -317.1546168986081
$arr1[5,2] = 0.6464407368884659
$arr1[6,0] = 'Season[T.Winter]'
$arr1[6,1] = -264.8053754024059
$arr1[6,2] = 0.6700531634783693
$arr1[7,0] = 'HHSize'
$arr1[7,1] = -391.8258384433015
$arr1[7,2] = 0.0552542286524214
$arr1[8,0] = 'Sex'
$arr1[8,1] = -441.0162766783388
$arr1[8,2] = 0.3445990408712415
$arr1[9,0] = 'Age'
$arr1[9,1] = -25.05240099010607
$arr1[9,2] = 0.238569259211287
$arr1[10,0] = 'DistSubcenter_res'
$arr1[10,1] = -1022.365436764895
$arr1[10,2] = 0.04721113951801786
$arr1[11,0] = 'DistCenter_res'
$arr1[11,1] = -591.1300974181202
$arr1[11,2] = 0.5450623612442389
$arr1[12,0] = 'UrbPopDensity_res'
$arr1[12,1] = 0.2657564068464622
$arr1[12,2] = 0.48986192362921
$arr1[13,0] = 'UrbBuildDensity_res'
$arr1[13,1] = -0.002355879623686632
$arr1[13,2] = 0.004369211380158078
$arr1[14,0] = 'IntersecDensity_res'
$arr1[14,1] = -130.1100682100447
$arr1[14,2] = 0.08350906483119774
$arr1[15,0] = 'street_length_res'
$arr1[15,1] = -148.0313809266345
$arr1[15,2] = 0.1535145896575459
$arr1[16,0] = 'LU_Comm_res'
$arr1[16,1] = 3547.479951120311
$arr1[16,2] = 0.6579015769250094
$arr1[17,0] = 'LU_UrbFab_res'
$arr1[17,1] = -2116.217396190945
$arr1[17,2] = 0.6901552554935746
$arr1[18,0] = 'bike_lane_share_res'
$arr1[18,1] = -15692.19564378664
$arr1[18,2] = 0.3327336821050153
$ws.Range("A2:C20").Value = $arr1

# Sheet 2: summ13615682 -> summ51175966
$ws = $wb.Worksheets.Item(2)
$ws.Name = 'summ51175966'
$arr2 = New-Object 'object[,]' 19,3
$arr2[0,0] = 'Intercept'
$arr2[0,1] = 28034.16086712027
$arr2[0,2] = 0.03072224058654197
$arr2[1,0] = 'Education[T.Primary/None]'
$arr2[1,1] = -777.5052919938994
$arr2[1,2] = 0.6599003427579078
$arr2[2,0] = 'Education[T.Secondary]'
$arr2[2,1] = -997.4447600548976
$arr2[2,2] = 0.4634905731596179
$arr2[3,0] = 'Education[T.University]'
$arr2[3,1] = -16.64810740911827
$arr2[3,2] = 0.9763956640335303
$arr2[4,0] = 'Season[T.Spring]'
$arr2[4,1] = 41.66331775925067
$arr2[4,2] = 0.9464040746221941
$arr2[5,0] = 'Season[T.Summer]'
$arr2[5,1] = -295.6877047680177
$arr2[5,2] = 0.668717698916002
$arr2[6,0] = 'Season[T.Winter]'
$arr2[6,1] = -625.3500573576781
$arr2[6,2] = 0.3145465675542303
$arr2[7,0] = 'HHSize'
$arr2[7,1] = -471.7040461434064
$arr2[7,2] = 0.02212696439364234
$arr2[8,0] = 'Sex'
$arr2[8,1] = -313.6241241481352
$arr2[8,2] = 0.5000890588596185
$arr2[9,0] = 'Age'
$arr2[9,1] = -14.8044405760428
$arr2[9,2] = 0.4960997053856014
$arr2[10,0] = 'DistSubcenter_res'
$arr2[10,1] = -512.9943037940202
$arr2[10,2] = 0.3133224194173985
$arr2[11,0] = 'DistCenter_res'
$arr2[11,1] = 47.74786377998154
$arr2[11,2] = 0.9599797993501418
$arr2[12,0] = 'UrbPopDensity_res'
$arr2[12,1] = 0.2448994247493877
$arr2[12,2] = 0.5182476068651067
$arr2[13,0] = 'UrbBuildDensity_res'
$arr2[13,1] = -0.001249096605059382
$arr2[13,2] = 0.1220608627252269
$arr2[14,0] = 'IntersecDensity_res'
$arr2[14,1] = -49.09252476974746
$arr2[14,2] = 0.5089682357966823
$arr2[15,0] = 'street_length_res'
$arr2[15,1] = -104.2769950834644
$arr2[15,2] = 0.315522381877953
$arr2[16,0] = 'LU_Comm_res'
$arr2[16,1] = 2769.036565294683
$arr2[16,2] = 0.7321294453562943
$arr2[17,0] = 'LU_UrbFab_res'
$arr2[17,1] = -5048.084915126261
$arr2[17,2] = 0.3381611323381204
$arr2[18,0] = 'bike_lane_share_res'
$arr2[18,1] = -23.23168942188568
$arr2[18,2] = 0.998838060053608
$ws.Range("A2:C20").Value = $arr2

# Sheet 3: summ13865359 -> summ51364523
$ws = $wb.Worksheets.Item(3)
$ws.Name = 'summ51364523'
$arr3 = New-Object 'object[,]' 19,3
$arr3[0,0] = 'Intercept'
$arr3[0,1] = 36206.18753036394
$arr3[0,2] = 0.007914025602784345
$arr3[1,0] = 'Education[T.Primary/None]'
$arr3[1,1] = -554.1027513821349
$arr3[1,2] = 0.7615292929464749
$arr3[2,0] = 'Education[T.Secondary]'
$arr3[2,1] = -613.263031964157
$arr3[2,2] = 0.6763785820129213
$arr3[3,0] = 'Education[T.University]'
$arr3[3,1] = 75.91799221022632
$arr3[3,2] = 0.8964875127585272
$arr3[4,0] = 'Season[T.Spring]'
$arr3[4,1] = -143.9243174352631
$arr3[4,2] = 0.8230156748628269
$arr3[5,0] = 'Season[T.Summer]'
$arr3[5,1] = -147.1755453356038
$arr3[5,2] = 0.8362939492267838
$arr3[6,0] = 'Season[T.Winter]'
$arr3[6,1] = -566.2830176753594
$arr3[6,2] = 0.3810968092981049
$arr3[7,0] = 'HHSize'
$arr3[7,1] = -338.1940111349891
$arr3[7,2] = 0.112674512365154
$arr3[8,0] = 'Sex'
$arr3[8,1] = 107.672362078677
$arr3[8,2] = 0.8231790671407003
$arr3[9,0] = 'Age'
$arr3[9,1] = -14.74057580254597
$arr3[9,2] = 0.5051208914947763
$arr3[10,0] = 'DistSubcenter_res'
$arr3[10,1] = -867.3017622729704
$arr3[10,2] = 0.10115067674528
$arr3[11,0] = 'DistCenter_res'
$arr3[11,1] = -677.292001172431
$arr3[11,2] = 0.4981350083728792
$arr3[12,0] = 'UrbPopDensity_res'
$arr3[12,1] = 0.3276227805018136
$arr3[12,2] = 0.3970985478756639
$arr3[13,0] = 'UrbBuildDensity_res'
$arr3[13,1] = -0.002047369049859044
$arr3[13,2] = 0.01601151980000711
$arr3[14,0] = 'IntersecDensity_res'
$arr3[14,1] = -97.50823195061707
$arr3[14,2] = 0.2126980236710896
$arr3[15,0] = 'street_length_res'
$arr3[15,1] = -114.1455760904581
$arr3[15,2] = 0.2844234922544477
$arr3[16,0] = 'LU_Comm_res'
$arr3[16,1] = 3167.944663625629
$arr3[16,2] = 0.7030171264074615
$arr3[17,0] = 'LU_UrbFab_res'
$arr3[17,1] = -5022.742582106788
$arr3[17,2] = 0.3436230359834124
$arr3[18,0] = 'bike_lane_share_res'
$arr3[18,1] = -11564.65842646863
$arr3[18,2] = 0.4843817886405616
$ws.Range("A2:C20").Value = $arr3

# Sheet 4: summ14131831 -> summ51550087
$ws = $wb.Worksheets.Item(4)
$ws.Name = 'summ51550087'
$arr4 = New-Object 'object[,]' 19,3
$arr4[0,0] = 'Intercept'
$arr4[0,1] = 35117.39568473565
$arr4[0,2] = 0.007420832143106108
$arr4[1,0] = 'Education[T.Primary/None]'
$arr4[1,1] = -621.8856731099704
$arr4[1,2] = 0.7423798730128341
$arr4[2,0] = 'Education[T.Secondary]'
$arr4[2,1] = -2105.940836901135
$arr4[2,2] = 0.1343967388436519
$arr4[3,0] = 'Education[T.University]'
$arr4[3,1] = -126.5715013043223
$arr4[3,2] = 0.8202665128543278
$arr4[4,0] = 'Season[T.Spring]'
$arr4[4,1] = -313.6161524333001
$arr4[4,2] = 0.6096495726280944
$arr4[5,0] = 'Season[T.Summer]'
$arr4[5,1] = 12.2619600764537
$arr4[5,2] = 0.9860304133951902
$arr4[6,0] = 'Season[T.Winter]'
$arr4[6,1] = -318.2109117962026
$arr4[6,2] = 0.6015813043468098
$arr4[7,0] = 'HHSize'
$arr4[7,1] = -357.98801973465
$arr4[7,2] = 0.08566865166373304
$arr4[8,0] = 'Sex'
$arr4[8,1] = 268.756429350738
$arr4[8,2] = 0.5619007669754557
$arr4[9,0] = 'Age'
$arr4[9,1] = -42.884297109933
$arr4[9,2] = 0.04457184279459201
$arr4[10,0] = 'DistSubcenter_res'
$arr4[10,1] = -675.495678875438
$arr4[10,2] = 0.189637000131452
$arr4[11,0] = 'DistCenter_res'
$arr4[11,1] = -1081.951968963463
$arr4[11,2] = 0.2678659545500088
$arr4[12,0] = 'UrbPopDensity_res'
$arr4[12,1] = 0.3311290921867955
$arr4[12,2] = 0.3825359983038202
$arr4[13,0] = 'UrbBuildDensity_res'
$arr4[13,1] = -0.001821021651705199
$arr4[13,2] = 0.02637975118903976
$arr4[14,0] = 'IntersecDensity_res'
$arr4[14,1] = -146.9766059066468
$arr4[14,2] = 0.04997657110075732
$arr4[15,0] = 'street_length_res'
$arr4[15,1] = -67.96702865776672
$arr4[15,2] = 0.5247194377392141
$arr4[16,0] = 'LU_Comm_res'
$arr4[16,1] = -1767.337947395502
$arr4[16,2] = 0.8296618328499918
$arr4[17,0] = 'LU_UrbFab_res'
$arr4[17,1] = -3906.404954206502
$arr4[17,2] = 0.4559927626291099
$arr4[18,0] = 'bike_lane_share_res'
$arr4[18,1] = -15981.41325996003
$arr4[18,2] = 0.3180183985429142
$ws.Range("A2:C20").Value = $arr4

# Sheet 5: summ14387273 -> summ51784883
$ws = $wb.Worksheets.Item(5)
$ws.Name = 'summ51784883'
$arr5 = New-Object 'object[,]' 19,3
$arr5[0,0] = 'Intercept'
$arr5[0,1] = 31485.29631222763
$arr5[0,2] = 0.01894055687531871
$arr5[1,0] = 'Education[T.Primary/None]'
$arr5[1,1] = -480.0607197986384
$arr5[1,2] = 0.7774479089402185
$arr5[2,0] = 'Education[T.Secondary]'
$arr5[2,1] = -1127.940946278559
$arr5[2,2] = 0.4007470197864719
$arr5[3,0] = 'Education[T.University]'
$arr5[3,1] = -65.67570827430978
$arr5[3,2] = 0.9070734445706918
$arr5[4,0] = 'Season[T.Spring]'
$arr5[4,1] = -621.528478558519
$arr5[4,2] = 0.3144336676050941
$arr5[5,0] = 'Season[T.Summer]'
$arr5[5,1] = -398.4504384620568
$arr5[5,2] = 0.5619094288247357
$arr5[6,0] = 'Season[T.Winter]'
$arr5[6,1] = -446.8686377468183
$arr5[6,2] = 0.4731325415968919
$arr5[7,0] = 'HHSize'
$arr5[7,1] = -353.1082285755958
$arr5[7,2] = 0.08398460594878981
$arr5[8,0] = 'Sex'
$arr5[8,1] = 355.8140539917268
$arr5[8,2] = 0.4462370810654753
$arr5[9,0] = 'Age'
$arr5[9,1] = 0.9664994038070702
$arr5[9,2] = 0.9637517889360581
$arr5[10,0] = 'DistSubcenter_res'
$arr5[10,1] = -769.0518260377867
$arr5[10,2] = 0.1367877558172311
$arr5[11,0] = 'DistCenter_res'
$arr5[11,1] = -767.1839123945821
$arr5[11,2] = 0.4331262856853829
$arr5[12,0] = 'UrbPopDensity_res'
$arr5[12,1] = 0.4090573411357458
$arr5[12,2] = 0.287270889325405
$arr5[13,0] = 'UrbBuildDensity_res'
$arr5[13,1] = -0.001852780593213048
$arr5[13,2] = 0.02573995075920915
$arr5[14,0] = 'IntersecDensity_res'
$arr5[14,1] = -110.088638358376
$arr5[14,2] = 0.1470460970013719
$arr5[15,0] = 'street_length_res'
$arr5[15,1] = -77.18428956036584
$arr5[15,2] = 0.4909121000592735
$arr5[16,0] = 'LU_Comm_res'
$arr5[16,1] = 1584.956109620803
$arr5[16,2] = 0.8512499471758503
$arr5[17,0] = 'LU_UrbFab_res'
$arr5[17,1] = -5451.100860436031
$arr5[17,2] = 0.3084948479528297
$arr5[18,0] = 'bike_lane_share_res'
$arr5[18,1] = -11981.00489873445
$arr5[18,2] = 0.4573926174247379
$ws.Range("A2:C20").Value = $arr5

# Sheet 6: summ14691676 -> summ51985209
$ws = $wb.Worksheets.Item(6)
$ws.Name = 'summ51985209'
$arr6 = New-Object 'object[,]' 19,3
$arr6[0,0] = 'Intercept'
$arr6[0,1] = 33338.78767554111
$arr6[0,2] = 0.01046592044793185
$arr6[1,0] = 'Education[T.Primary/None]'
$arr6[1,1] = -302.6049730480615
$arr6[1,2] = 0.8562490718635984
$arr6[2,0] = 'Education[T.Secondary]'
$arr6[2,1] = -1397.091056309723
$arr6[2,2] = 0.2932118512322415
$arr6[3,0] = 'Education[T.University]'
$arr6[3,1] = 180.1946333506608
$arr6[3,2] = 0.7498701464396821
$arr6[4,0] = 'Season[T.Spring]'
$arr6[4,1] = -229.0089329747622
$arr6[4,2] = 0.711757919106367
$arr6[5,0] = 'Season[T.Summer]'
$arr6[5,1] = -298.2203226120951
$arr6[5,2] = 0.671416197929129
$arr6[6,0] = 'Season[T.Winter]'
$arr6[6,1] = -353.2655978264316
$arr6[6,2] = 0.5690528865112352
$arr6[7,0] = 'HHSize'
$arr6[7,1] = -286.4809925302326
$arr6[7,2] = 0.1640178171701589
$arr6[8,0] = 'Sex'
$arr6[8,1] = 167.3468976340443
$arr6[8,2] = 0.7213579781380011
$arr6[9,0] = 'Age'
$arr6[9,1] = -24.24736158127233
$arr6[9,2] = 0.254706930794961
$arr6[10,0] = 'DistSubcenter_res'
$arr6[10,1] = -857.6333719849908
$arr6[10,2] = 0.09234316531512636
$arr6[11,0] = 'DistCenter_res'
$arr6[11,1] = -655.6739792682001
$arr6[11,2] = 0.4828750697010978
$arr6[12,0] = 'UrbPopDensity_res'
$arr6[12,1] = 0.1057682516251894
$arr6[12,2] = 0.7857991680614608
$arr6[13,0] = 'UrbBuildDensity_res'
$arr6[13,1] = -0.001742222926652802
$arr6[13,2] = 0.03095743392550081
$arr6[14,0] = 'IntersecDensity_res'
$arr6[14,1] = -89.97541167871347
$arr6[14,2] = 0.2226922639613086
$arr6[15,0] = 'street_length_res'
$arr6[15,1] = -94.05117103297592
$arr6[15,2] = 0.3671057043603873
$arr6[16,0] = 'LU_Comm_res'
$arr6[16,1] = 958.1781028886498
$arr6[16,2] = 0.9039881817273381
$arr6[17,0] = 'LU_UrbFab_res'
$arr6[17,1] = -5604.652524106135
$arr6[17,2] = 0.301760588618181
$arr6[18,0] = 'bike_lane_share_res'
$arr6[18,1] = 502.6735486894768
$arr6[18,2] = 0.9748900914308682
$ws.Range("A2:C20").Value = $arr6

# Sheet 7: summ14954215 -> summ52177872
$ws = $wb.Worksheets.Item(7)
$ws.Name = 'summ52177872'
$arr7 = New-Object 'object[,]' 19,3
$arr7[0,0] = 'Intercept'
$arr7[0,1] = 31583.87077206546
$arr7[0,2] = 0.01932333410543951
$arr7[1,0] = 'Education[T.Primary/None]'
$arr7[1,1] = -573.3881976337714
$arr7[1,2] = 0.7375933997718889
$arr7[2,0] = 'Education[T.Secondary]'
$arr7[2,1] = -1319.153264391747
$arr7[2,2] = 0.3420727717606817
$arr7[3,0] = 'Education[T.University]'
$arr7[3,1] = 272.3607020428479
$arr7[3,2] = 0.6415937627837592
$arr7[4,0] = 'Season[T.Spring]'
$arr7[4,1] = -74.22293191276542
$arr7[4,2] = 0.9075940948893193
$arr7[5,0] = 'Season[T.Summer]'
$arr7[5,1] = -214.297080489913
$arr7[5,2] = 0.7625694856809203
$arr7[6,0] = 'Season[T.Winter]'
$arr7[6,1] = -434.2591436103848
$arr7[6,2] = 0.4958302739968437
$arr7[7,0] = 'HHSize'
$arr7[7,1] = -328.3707426587562
$arr7[7,2] = 0.1206248150320475
$arr7[8,0] = 'Sex'
$arr7[8,1] = -32.26410534765489
$arr7[8,2] = 0.9462492731176847
$arr7[9,0] = 'Age'
$arr7[9,1] = -21.17348786636105
$arr7[9,2] = 0.3367395388441242
$arr7[10,0] = 'DistSubcenter_res'
$arr7[10,1] = -562.1628644296762
$arr7[10,2] = 0.2846014473163202
$arr7[11,0] = 'DistCenter_res'
$arr7[11,1] = -282.7553113457642
$arr7[11,2] = 0.7784740700820076
$arr7[12,0] = 'UrbPopDensity_res'
$arr7[12,1] = 0.3580184530535697
$arr7[12,2] = 0.3493079972681695
$arr7[13,0] = 'UrbBuildDensity_res'
$arr7[13,1] = -0.001620780625029019
$arr7[13,2] = 0.05468168301841373
$arr7[14,0] = 'IntersecDensity_res'
$arr7[14,1] = -94.35258500726565
$arr7[14,2] = 0.2286201696835723
$arr7[15,0] = 'street_length_res'
$arr7[15,1] = -107.4650013690638
$arr7[15,2] = 0.3240160619087121
$arr7[16,0] = 'LU_Comm_res'
$arr7[16,1] = 3519.710552020207
$arr7[16,2] = 0.6783132838095687
$arr7[17,0] = 'LU_UrbFab_res'
$arr7[17,1] = -5086.574265493851
$arr7[17,2] = 0.3290603145791194
$arr7[18,0] = 'bike_lane_share_res'
$arr7[18,1] = -30.32200403416755
$arr7[18,2] = 0.998563550037923
$ws.Range("A2:C20").Value = $arr7

# Sheet 8: summ15218137 -> summ52368543
$ws = $wb.Worksheets.Item(8)
$ws.Name = 'summ52368543'
$arr8 = New-Object 'object[,]' 19,3
$arr8[0,0] = 'Intercept'
$arr8[0,1] = 34297.92942277397
$arr8[0,2] = 0.01087832623249684
$arr8[1,0] = 'Education[T.Primary/None]'
$arr8[1,1] = -588.799302858888
$arr8[1,2] = 0.734092870676026
$arr8[2,0] = 'Education[T.Secondary]'
$arr8[2,1] = -1665.682121930077
$arr8[2,2] = 0.2473334064528379
$arr8[3,0] = 'Education[T.University]'
$arr8[3,1] = 279.7074267956456
$arr8[3,2] = 0.6228102093245094
$arr8[4,0] = 'Season[T.Spring]'
$arr8[4,1] = -64.26134972075852
$arr8[4,2] = 0.9189423107365651
$arr8[5,0] = 'Season[T.Summer]'
$arr8[5,1] = -144.4659789011504
$arr8[5,2] = 0.8375991849688947
$arr8[6,0] = 'Season[T.Winter]'
$arr8[6,1] = -626.3684197951649
$arr8[6,2] = 0.3252136856456326
$arr8[7,0] = 'HHSize'
$arr8[7,1] = -384.5872095213876
$arr8[7,2] = 0.06672112414436036
$arr8[8,0] = 'Sex'
$arr8[8,1] = -299.3323210744891
$arr8[8,2] = 0.5316537678163327
$arr8[9,0] = 'Age'
$arr8[9,1] = -8.904790831589773
$arr8[9,2] = 0.6816792148773834
$arr8[10,0] = 'DistSubcenter_res'
$arr8[10,1] = -850.6639180409115
$arr8[10,2] = 0.09979828561149769
$arr8[11,0] = 'DistCenter_res'
$arr8[11,1] = -629.0286800611593
$arr8[11,2] = 0.5304668696653838
$arr8[12,0] = 'UrbPopDensity_res'
$arr8[12,1] = 0.3403492344101622
$arr8[12,2] = 0.3964855998712502
$arr8[13,0] = 'UrbBuildDensity_res'
$arr8[13,1] = -0.001913194840080912
$arr8[13,2] = 0.02196236836647992
$arr8[14,0] = 'IntersecDensity_res'
$arr8[14,1] = -101.0065358715659
$arr8[14,2] = 0.190977114224281
$arr8[15,0] = 'street_length_res'
$arr8[15,1] = -103.5490709120535
$arr8[15,2] = 0.3537760986777512
$arr8[16,0] = 'LU_Comm_res'
$arr8[16,1] = 2035.972489851726
$arr8[16,2] = 0.8147227428606683
$arr8[17,0] = 'LU_UrbFab_res'
$arr8[17,1] = -4722.477601124961
$arr8[17,2] = 0.3909985085041152
$arr8[18,0] = 'bike_lane_share_res'
$arr8[18,1] = -8533.174546992872
$arr8[18,2] = 0.6013681371653999
$ws.Range("A2:C20").Value = $arr8

# Sheet 9: summ15489248 -> summ52552994
$ws = $wb.Worksheets.Item(9)
$ws.Name = 'summ52552994'
$arr9 = New-Object 'object[,]' 19,3
$arr9[0,0] = 'Intercept'
$arr9[0,1] = 39089.87081477432
$arr9[0,2] = 0.003137839463808443
$arr9[1,0] = 'Education[T.Primary/None]'
$arr9[1,1] = 442.795826782806
$arr9[1,2] = 0.8102507407576129
$arr9[2,0] = 'Education[T.Secondary]'
$arr9[2,1] = -984.8851320052667
$arr9[2,2] = 0.4747120671133502
$arr9[3,0] = 'Education[T.University]'
$arr9[3,1] = 457.6253036918075
$arr9[3,2] = 0.4236914077545172
$arr9[4,0] = 'Season[T.Spring]'
$arr9[4,1] = 111.8835456173952
$arr9[4,2] = 0.8603581778329092
$arr9[5,0] = 'Season[T.Summer]'
$arr9[5,1] = -54.95769230061859
$arr9[5,2] = 0.9370403207906108
$arr9[6,0] = 'Season[T.Winter]'
$arr9[6,1] = -536.4652030017345
$arr9[6,2] = 0.3987543084420112
$arr9[7,0] = 'HHSize'
$arr9[7,1] = -377.6279543245055
$arr9[7,2] = 0.06986003620348478
$arr9[8,0] = 'Sex'
$arr9[8,1] = 83.89870558571351
$arr9[8,2] = 0.8600159334964161
$arr9[9,0] = 'Age'
$arr9[9,1] = -16.08856614343748
$arr9[9,2] = 0.4602265867531768
$arr9[10,0] = 'DistSubcenter_res'
$arr9[10,1] = -1038.010302170077
$arr9[10,2] = 0.04552034383583096
$arr9[11,0] = 'DistCenter_res'
$arr9[11,1] = -1273.53753206139
$arr9[11,2] = 0.1943410758380871
$arr9[12,0] = 'UrbPopDensity_res'
$arr9[12,1] = 0.200079331940524
$arr9[12,2] = 0.6014914150467159
$arr9[13,0] = 'UrbBuildDensity_res'
$arr9[13,1] = -0.002171302448566961
$arr9[13,2] = 0.00876631488465664
$arr9[14,0] = 'IntersecDensity_res'
$arr9[14,1] = -154.0627356646911
$arr9[14,2] = 0.04282534977593552
$arr9[15,0] = 'street_length_res'
$arr9[15,1] = -85.26733097342967
$arr9[15,2] = 0.4185179284703289
$arr9[16,0] = 'LU_Comm_res'
$arr9[16,1] = -2428.373677123508
$arr9[16,2] = 0.7693128022978585
$arr9[17,0] = 'LU_UrbFab_res'
$arr9[17,1] = -3267.239945671324
$arr9[17,2] = 0.534169714647835
$arr9[18,0] = 'bike_lane_share_res'
$arr9[18,1] = -24718.32633887539
$arr9[18,2] = 0.1304119570785365
$ws.Range("A2:C20").Value = $arr9

